$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 5941
$ws.Range("I4").Value = 6352.7144
$ws.Range("J4").Value = 4500
$ws.Range("K4").Value = 6352.7144
$ws.Range("L4").Value = 4500
$ws.Range("M4").Value = -6238.7144
$ws.Range("N4").Value = -4728

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1597.5
$ws.Range("I62").Value = 1597.5
$ws.Range("K62").Value = 1597.5
$ws.Range("M62").Value = -973.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 1597.5
$ws.Range("I65").Value = 1597.5
$ws.Range("K65").Value = 7987.5
$ws.Range("M65").Value = -4867.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2633.3333
$ws.Range("J70").Value = 2500
$ws.Range("L70").Value = 7500
$ws.Range("N70").Value = -8040

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 2633.3333
$ws.Range("J73").Value = 2500
$ws.Range("L73").Value = 7500
$ws.Range("N73").Value = -9372

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 733.4
$ws.Range("I107").Value = 733.4
$ws.Range("K107").Value = 733.4
$ws.Range("M107").Value = 1186.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7050.25
$ws.Range("I116").Value = 6483.1665
$ws.Range("J116").Value = 8751.5
$ws.Range("K116").Value = 6483.1665
$ws.Range("L116").Value = 8751.5
$ws.Range("M116").Value = -3041.1665
$ws.Range("N116").Value = -15635.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 72000
$ws.Range("J123").Value = 72000
$ws.Range("L123").Value = 72000
$ws.Range("N123").Value = -81800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1247.6
$ws.Range("J125").Value = 199.5
$ws.Range("L125").Value = 1795.5
$ws.Range("N125").Value = -6715.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 69999
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("N133").Value = -80119

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2854.4
$ws.Range("I63").Value = 2854.4
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 2854.4
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2168.4
$ws.Range("N63").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 2854.4
$ws.Range("I66").Value = 2854.4
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 14272
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -10840
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6392.2896
$ws.Range("I132").Value = 4288.0347
$ws.Range("K132").Value = 12864.1041
$ws.Range("M132").Value = -10334.1041

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 34322.168
$ws.Range("J100").Value = 34322.168
$ws.Range("L100").Value = 34322.168
$ws.Range("N100").Value = -36486.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 41161.37
$ws.Range("I134").Value = 1559.6666
$ws.Range("J134").Value = 90663.5
$ws.Range("K134").Value = 4678.9998
$ws.Range("L134").Value = 271990.5
$ws.Range("M134").Value = -2143.9998
$ws.Range("N134").Value = -277060.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 480495.4
$ws.Range("I31").Value = 3353.875
$ws.Range("J31").Value = 938551.25
$ws.Range("K31").Value = 3353.875
$ws.Range("L31").Value = 938551.25
$ws.Range("M31").Value = -3058.875
$ws.Range("N31").Value = -939141.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 480495.4
$ws.Range("I34").Value = 3353.875
$ws.Range("J34").Value = 938551.25
$ws.Range("K34").Value = 3353.875
$ws.Range("L34").Value = 938551.25
$ws.Range("M34").Value = -3151.875
$ws.Range("N34").Value = -938955.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 958.03125
$ws.Range("I58").Value = 656.0714
$ws.Range("K58").Value = 656.0714
$ws.Range("M58").Value = -453.0714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4249
$ws.Range("I99").Value = 3687
$ws.Range("J99").Value = 4998.3335
$ws.Range("K99").Value = 3687
$ws.Range("L99").Value = 4998.3335
$ws.Range("M99").Value = -2189
$ws.Range("N99").Value = -7994.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 49500
$ws.Range("J104").Value = 49500
$ws.Range("L104").Value = 49500
$ws.Range("N104").Value = -54742

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2439.2856
$ws.Range("I105").Value = 1942.75
$ws.Range("J105").Value = 3101.3333
$ws.Range("K105").Value = 1942.75
$ws.Range("L105").Value = 3101.3333
$ws.Range("M105").Value = -195.75
$ws.Range("N105").Value = -6595.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4249
$ws.Range("I126").Value = 3687
$ws.Range("J126").Value = 4998.3335
$ws.Range("K126").Value = 11061
$ws.Range("L126").Value = 14995.0005
$ws.Range("M126").Value = -8591
$ws.Range("N126").Value = -19935.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4323.909
$ws.Range("I132").Value = 4283.222
$ws.Range("K132").Value = 12849.666
$ws.Range("M132").Value = -10319.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1436633.6
$ws.Range("I134").Value = 5001994.5
$ws.Range("K134").Value = 15005983.5
$ws.Range("M134").Value = -15003448.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 958.03125
$ws.Range("I136").Value = 656.0714
$ws.Range("K136").Value = 1968.2142
$ws.Range("M136").Value = 581.7857999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4252.0835
$ws.Range("I80").Value = 2946.5
$ws.Range("K80").Value = 8839.5
$ws.Range("M80").Value = -7903.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 4252.0835
$ws.Range("I83").Value = 2946.5
$ws.Range("K83").Value = 26518.5
$ws.Range("M83").Value = -21838.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 2357.5625
$ws.Range("I139").Value = 2740.6667
$ws.Range("J139").Value = 2269.1538
$ws.Range("K139").Value = 8222.000100000001
$ws.Range("L139").Value = 6807.4614
$ws.Range("M139").Value = -3082.000100000001
$ws.Range("N139").Value = -17087.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 217379.64
$ws.Range("I140").Value = 275560.72
$ws.Range("K140").Value = 826682.1599999999
$ws.Range("M140").Value = -821502.1599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2112.7144
$ws.Range("I122").Value = 1758
$ws.Range("K122").Value = 5274
$ws.Range("M122").Value = -2824

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4073.5
$ws.Range("I126").Value = 3541.5
$ws.Range("J126").Value = 4339.5
$ws.Range("K126").Value = 10624.5
$ws.Range("L126").Value = 13018.5
$ws.Range("M126").Value = -8154.5
$ws.Range("N126").Value = -17958.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 62666.668
$ws.Range("J139").Value = 62666.668
$ws.Range("L139").Value = 62666.668
$ws.Range("N139").Value = -72946.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1949.2778
$ws.Range("I40").Value = 1139.1333
$ws.Range("J40").Value = 6000
$ws.Range("K40").Value = 1139.1333
$ws.Range("L40").Value = 6000
$ws.Range("M40").Value = -1003.1333
$ws.Range("N40").Value = -6272

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 909.5833
$ws.Range("I82").Value = 576.6667
$ws.Range("J82").Value = 1242.5
$ws.Range("K82").Value = 576.6667
$ws.Range("L82").Value = 1242.5
$ws.Range("M82").Value = -215.6667
$ws.Range("N82").Value = -1964.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 909.5833
$ws.Range("I85").Value = 576.6667
$ws.Range("J85").Value = 1242.5
$ws.Range("K85").Value = 576.6667
$ws.Range("L85").Value = 1242.5
$ws.Range("M85").Value = 671.3333
$ws.Range("N85").Value = -3738.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 99600.63
$ws.Range("I132").Value = 78645.92
$ws.Range("J132").Value = 129868.555
$ws.Range("K132").Value = 235937.76
$ws.Range("L132").Value = 389605.665
$ws.Range("M132").Value = -233407.76
$ws.Range("N132").Value = -394665.665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 120772.4
$ws.Range("I136").Value = 8531.286
$ws.Range("J136").Value = 382668.34
$ws.Range("K136").Value = 25593.858
$ws.Range("L136").Value = 1148005.02
$ws.Range("M136").Value = -23043.858
$ws.Range("N136").Value = -1153105.02

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 50002040
$ws.Range("I107").Value = 62501940
$ws.Range("J107").Value = 2451.5
$ws.Range("K107").Value = 187505820
$ws.Range("L107").Value = 7354.5
$ws.Range("M107").Value = -187503900
$ws.Range("N107").Value = -11194.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3689.7144
$ws.Range("I126").Value = 2115.6
$ws.Range("K126").Value = 6346.799999999999
$ws.Range("M126").Value = -3876.799999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8709.5
$ws.Range("I132").Value = 2120.4119
$ws.Range("J132").Value = 14605
$ws.Range("K132").Value = 6361.2357
$ws.Range("L132").Value = 43815
$ws.Range("M132").Value = -3831.2357
$ws.Range("N132").Value = -48875
